$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename "Sheet2" -> "BOB"
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Name = "BOB"

# ------------------------------------------------------------------
# 2. Populate the BOB sheet with the FD (Fixed Deposit) report data
# ------------------------------------------------------------------
$rupee = [char]8377

# Header row 1 -- all quote-prefixed ('s="1"' / quotePrefix style) except
# FDTenureDays (column D), which is written as a plain string.
$ws.Cells.Item(1, 1).Value = "'FDDepositType"
$ws.Cells.Item(1, 2).Value = "'FDAmount"
$ws.Cells.Item(1, 3).Value = "'FDTenure"
$ws.Cells.Item(1, 4).Value = "FDTenureDays"
$ws.Cells.Item(1, 5).Value = "'FDRateOfInterest"
$ws.Cells.Item(1, 6).Value = "'FDMaturityDate"
$ws.Cells.Item(1, 7).Value = "'FDMaturityValue"
$ws.Cells.Item(1, 8).Value = "'AggregateInterestAmount"
$ws.Cells.Item(1, 9).Value = "'InterestPerQuarter"

# Data rows 2-11. Text values that look numeric must be forced to text
# (apostrophe prefix) then reset back to the "Normal" style so they keep
# type "s" (shared string) but do not pick up the quote-prefix style.
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = "Quarterly Payout"

    $ws.Cells.Item($r, 2).Value = "'100000"
    $ws.Cells.Item($r, 2).Style = "Normal"

    $ws.Cells.Item($r, 3).Value = "5 Years: 0 Months : 0 Days"

    $ws.Cells.Item($r, 4).Value = "'1826"
    $ws.Cells.Item($r, 4).Style = "Normal"

    $ws.Cells.Item($r, 5).Value = "'6.5"
    $ws.Cells.Item($r, 5).Style = "Normal"

    if ($r -le 6) {
        $ws.Cells.Item($r, 6).Value = "7 Mar 2029"
    } else {
        $ws.Cells.Item($r, 6).Value = "8 Mar 2029"
    }

    $ws.Cells.Item($r, 7).Value = "'100000"
    $ws.Cells.Item($r, 7).Style = "Normal"

    $ws.Cells.Item($r, 8).Value = ("'" + $rupee + "32,500")
    $ws.Cells.Item($r, 8).Style = "Normal"

    $ws.Cells.Item($r, 9).Value = "'1,625"
    $ws.Cells.Item($r, 9).Style = "Normal"
}

# ------------------------------------------------------------------
# 3. Column widths for the BOB sheet (values compensate for the
#    engine's internal character-width padding so the saved <col>
#    width matches the target as closely as possible).
# ------------------------------------------------------------------
$ws.Range("B1:E1").ColumnWidth = 15.166666666666666
$ws.Range("F1").ColumnWidth = 19.166666666666668
$ws.Range("G1").ColumnWidth = 17.307291666666668
$ws.Range("H1").ColumnWidth = 22.307291666666668

# ------------------------------------------------------------------
# 4. View state: make BOB the active/visible sheet, scroll so column C
#    is the left-most visible column, and select G3.
# ------------------------------------------------------------------
$ws.Activate()
$null = $ws.Range("G3").Select()
$excel.ActiveWindow.ScrollColumn = 3
